# Apply the log.xlsx edit:
#  - Clear the stray empty inline-string cells C11:E11 (they become truly
#    empty / absent cells, shrinking row 11 back down).
#  - Append a new row 12 with Olga's second entry (age 29, Drama/English/2001)
#    including the pandas-printed recommendation block in F12 and "Yes" in G12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: drop the leftover empty inline strings in C11, D11, E11.
$ws.Range("C11:E11").ClearContents()

# Row 12: new record for Olga. New rows otherwise pick up the column's
# default style (style index 1, wrapText) from <cols>; reset to "Normal"
# so the written cells carry no explicit per-cell style, matching the rest
# of the un-styled data rows (6-11).
$ws.Range("A12:G12").Style = "Normal"

$ws.Range("A12").Value = "Olga"
$ws.Range("B12").Value = 29
$ws.Range("C12").Value = "Drama"
$ws.Range("D12").Value = "English"

# Force E12 to be stored as text "2001" rather than being auto-coerced to a
# number: temporarily format as Text ("@") so the literal string sticks,
# then reset the style again so the cell carries no explicit style.
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2001"
$ws.Range("E12").Style = "Normal"

$f12 = @"
                 original_title  year  genre language      budget  duration usa_gross_income  reviews_from_users
31327            Ricochet River  2001  Drama  English         NaN       112              NaN                 7.0
31335         The Shipping News  2001  Drama  English  $ 38000000       111       $ 11434216               210.0
34751      The Invisible Circus  2001  Drama  English         NaN        93          $ 77578                49.0
35164  A Day in Black and White  2001  Drama  English         NaN        80              NaN                 4.0
36026          The Rising Place  2001  Drama  English         NaN        93           $ 8360                 2.0
36414       The Sleepy Time Gal  2001  Drama  English         NaN       108              NaN                14.0
36460    Diary of a City Priest  2001  Drama  English         NaN        77              NaN                 4.0
36705           Sensual Friends  2001  Drama  English         NaN        93              NaN                 6.0
36899    Goodbye Charlie Bright  2001  Drama  English         NaN        87              NaN                24.0
37077         Borderline Normal  2001  Drama  English         NaN        85              NaN                 8.0
"@

$ws.Range("F12").Value = $f12.TrimEnd("`r", "`n")
$ws.Range("G12").Value = "Yes"
